$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 62-69: "Tipo_Horario" (column N) changes from "Frecuencia" to "Flexible"
$ws.Range("N62:N69").Value = "Flexible"

# Column O ("Duracion_Trayecto_Min") picks up the minute values that used to
# live in column P ("Frecuencia_Min"), expressed as a fraction of a day
# (minutes / 1440) so the stored double matches Excel's time serialization.
$ws.Range("O62").Value = 8/1440
$ws.Range("O63").Value = 8/1440
$ws.Range("O64").Value = 12/1440
$ws.Range("O65").Value = 12/1440
$ws.Range("O66").Value = 25/1440
$ws.Range("O67").Value = 25/1440
$ws.Range("O68").Value = 35/1440
$ws.Range("O69").Value = 35/1440

# Column P ("Frecuencia_Min") is no longer used for these rows; fully clear
# the cells (not just their contents) so they disappear from the sheet XML.
$ws.Range("P62:P69").Clear()

# The now-narrower column P content lets the column shrink back down.
$ws.Columns.Item(16).ColumnWidth = 14.15

# Update the view: selection moves to G2 and the window scrolls back to the
# top of the sheet (no more frozen/forced topLeftCell on row 43).
$ws.Range("G2").Select()
